$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75; existing rows 75:98 shift down to 76:99
$ws.Rows.Item(75).Insert()

# Fill in the new row 75 with data (mirrors the other rows in the block)
$ws.Range("A75").Value = 11
$ws.Range("B75").Value = "Vega Monumental Concepción"
$ws.Range("C75").Value = "Bíobío"
$ws.Range("D75").Value = 44889
$ws.Range("E75").Value = 8
$ws.Range("F75").Value = "Fruta"
$ws.Range("G75").Value = 100101
$ws.Range("H75").Value = "Berries"
$ws.Range("I75").Value = 100101001
$ws.Range("J75").Value = "Arándano (blue)"
$ws.Range("K75").Value = "Sin especificar"
$ws.Range("L75").Value = "Primera"
$ws.Range("M75").Value = 200
$ws.Range("N75").Value = 4000
$ws.Range("O75").Value = 4500
$ws.Range("P75").Value = 4250
$ws.Range("Q75").Value = "$/bandeja 2 kilos"
$ws.Range("R75").Value = "Región de Ñuble"
$ws.Range("S75").Value = 2125
$ws.Range("T75").Value = 2

# Match the date cell style used by the other "Fecha" column cells (D column)
$ws.Range("D75").NumberFormat = $ws.Range("D76").NumberFormat
